$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new "Price" text would otherwise be auto-parsed by Excel as a
# number (losing formatting like trailing zeros, e.g. "0.790" -> 0.79). Force
# them to Text format first so the literal string is preserved, matching the
# source data which stores these as plain text.
$textCells = @("D5", "D6", "D7", "D9", "D10", "D11", "D12", "D16", "D18", "D20", "D21", "D22", "D23", "D24", "D26", "D27", "D29", "D32", "D33", "D34", "D35", "D38", "D39", "D40", "D41", "D42", "D44", "D46", "D47", "D48", "D49", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "46.781.19"
$ws.Range("E2").Value = "  -0.30%  "

$ws.Range("D3").Value = "2.262.63"
$ws.Range("E3").Value = "  -3.69%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").Value = "298.23"
$ws.Range("E5").Value = "  -2.86%  "

$ws.Range("D6").Value = "97.89"
$ws.Range("E6").Value = "  -0.36%  "

$ws.Range("D7").Value = "0.575"
$ws.Range("E7").Value = "  -0.54%  "

$ws.Range("E8").Value = "  +0.04%  "

$ws.Range("D9").Value = "0.501"
$ws.Range("E9").Value = "  -6.90%  "

$ws.Range("D10").Value = "34.64"
$ws.Range("E10").Value = "  -3.29%  "

$ws.Range("D11").Value = "0.0791"
$ws.Range("E11").Value = "  -1.93%  "

$ws.Range("D12").Value = "6.98"
$ws.Range("E12").Value = "  -6.30%  "

$ws.Range("E13").Value = "  -1.86%  "

$ws.Range("D14").Value = "2.609.26"
$ws.Range("E14").Value = "  -3.61%  "

$ws.Range("D15").Value = "2.267.13"
$ws.Range("E15").Value = "  -3.51%  "

$ws.Range("D16").Value = "13.56"
$ws.Range("E16").Value = "  -4.85%  "

$ws.Range("D17").Value = "46.780.79"
$ws.Range("E17").Value = "  -0.03%  "

$ws.Range("D18").Value = "0.790"
$ws.Range("E18").Value = "  -5.10%  "

$ws.Range("D19").Value = "0.0₃0963"
$ws.Range("E19").Value = "  +1.28%  "

$ws.Range("D20").Value = "12.36"
$ws.Range("E20").Value = "  -10.38%  "

$ws.Range("D21").Value = "5.76"
$ws.Range("E21").Value = "  -6.89%  "

$ws.Range("D22").Value = "65.68"
$ws.Range("E22").Value = "  -1.83%  "

$ws.Range("D23").Value = "244.41"
$ws.Range("E23").Value = "  -0.41%  "

$ws.Range("D24").Value = "2.77"
$ws.Range("E24").Value = "  -7.42%  "

$ws.Range("E25").Value = "  +0.87%  "

$ws.Range("D26").Value = "1.84"
$ws.Range("E26").Value = "  -8.05%  "

$ws.Range("D27").Value = "40.96"
$ws.Range("E27").Value = "  -2.43%  "

$ws.Range("E28").Value = "  -3.55%  "

$ws.Range("D29").Value = "9.48"
$ws.Range("E29").Value = "  -4.26%  "

$ws.Range("E30").Value = "  -0.77%  "

$ws.Range("E31").Value = "  +7.18%  "

$ws.Range("D32").Value = "3.30"
$ws.Range("E32").Value = "  +4.17%  "

$ws.Range("D33").Value = "144.81"
$ws.Range("E33").Value = "  -5.06%  "

$ws.Range("D34").Value = "5.27"
$ws.Range("E34").Value = "  -8.92%  "

$ws.Range("D35").Value = "0.0760"
$ws.Range("E35").Value = "  -6.70%  "

$ws.Range("E36").Value = "  +0.97%  "

$ws.Range("E37").Value = "  -2.94%  "

$ws.Range("D38").Value = "15.29"
$ws.Range("E38").Value = "  +10.63%  "

$ws.Range("D39").Value = "1.64"
$ws.Range("E39").Value = "  -10.46%  "

$ws.Range("D40").Value = "3.79"
$ws.Range("E40").Value = "  -6.80%  "

$ws.Range("D41").Value = "0.0292"
$ws.Range("E41").Value = "  -7.63%  "

$ws.Range("D42").Value = "3.05"
$ws.Range("E42").Value = "  -11.28%  "

$ws.Range("E43").Value = "  +0.05%  "

$ws.Range("D44").Value = "92.94"
$ws.Range("E44").Value = "  +14.71%  "

$ws.Range("D45").Value = "1.779.50"
$ws.Range("E45").Value = "  -5.55%  "

$ws.Range("D46").Value = "1.86"
$ws.Range("E46").Value = "  -6.88%  "

$ws.Range("D47").Value = "69.76"
$ws.Range("E47").Value = "  -6.23%  "

$ws.Range("D48").Value = "0.182"
$ws.Range("E48").Value = "  -7.88%  "

$ws.Range("D49").Value = "4.77"
$ws.Range("E49").Value = "  -3.47%  "

$ws.Range("D51").Value = "93.46"
$ws.Range("E51").Value = "  -5.64%  "

# Row 50: coin replaced (FraxShare -> RocketPoolETH) along with its link/price/volume
$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D50").Value = "2.485.86"
$ws.Range("E50").Value = "  -3.84%  "
